$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the added "OpenMP" speed-up columns
$ws.Range("J1").Value = "OpenMP"
$ws.Range("K1").Value = "SpeedUṕ"

# Row 2 (Serial C): OpenMP time + speed-up formula
$ws.Range("J2").Value = 768.09
$ws.Range("K2").Formula = "=B2/J2"
$ws.Range("K2").Font.Bold = $false

# Row 3 (Serial C++): OpenMP time + speed-up formula
$ws.Range("J3").Value = 1423.99
$ws.Range("K3").Formula = "=B3/J3"
$ws.Range("K3").Font.Bold = $true

# G3 loses its bold formatting (style changes from bold 0.00 to plain 0.00)
$ws.Range("G3").Font.Bold = $false

# Update the active selection to G3
$ws.Range("G3").Select()
